# Nexial base macro test fixture update:
# Add new "image" command ocr(source,saveVar) / ocr(image,saveVar)
# and new "tn.5250" (profile-based) command group: close/open/saveText/typeKeys/updateScreenFields
#
# This touches the "#system" worksheet, which holds one "named column" of
# command strings per command-group (image, target, web, webalert, webcookie,
# ws, ws.async, xml, ...). A brand new column (tn.5250) is inserted right
# before the old "web" column, which pushes web/webalert/webcookie/ws/
# ws.async/xml one column to the right. The "target" column (A) gets one
# new row (tn.5250) inserted before "web", and the "image" column (K) gets
# one new row (ocr(image,saveVar)) inserted before "resize(...)".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# --- 1) Insert a brand-new column at Z, shifting the old web..xml columns
#        (Z..AE) one column to the right (AA..AF). ---------------------
$ws.Columns("Z").Insert()

# Populate the new "tn.5250" column (Z2:Z6)
$ws.Range("Z1").Value = "tn.5250"
$ws.Range("Z2").Value = "close(profile)"
$ws.Range("Z3").Value = "open(profile)"
$ws.Range("Z4").Value = "saveText(profile,var)"
$ws.Range("Z5").Value = "typeKeys(profile,keystrokes)"
$ws.Range("Z6").Value = "updateScreenFields(profile)"

# --- 2) Add new "ocr(image,saveVar)" row into the "image" column (K),
#        inserted ahead of "resize(image,width,height,saveTo)", and rename
#        the existing "colorbit(source,...)" entry to "colorbit(image,...)".
$ws.Range("K8").Value = "saveDiff(var,baseline,actual)"
$ws.Range("K7").Value = "resize(image,width,height,saveTo)"
$ws.Range("K6").Value = "ocr(image,saveVar)"
$ws.Range("K2").Value = "colorbit(image,bit,saveTo)"

# --- 3) Add new "tn.5250" row into the "target" column (A), inserted
#        ahead of "web". -------------------------------------------------
$ws.Range("A32").Value = "xml"
$ws.Range("A31").Value = "ws.async"
$ws.Range("A30").Value = "ws"
$ws.Range("A29").Value = "webcookie"
$ws.Range("A28").Value = "webalert"
$ws.Range("A27").Value = "web"
$ws.Range("A26").Value = "tn.5250"

# --- 4) Fix up the named ranges so they describe the new layout. --------
$wb.Names.Item("image").RefersTo = "='#system'!`$K`$2:`$K`$8"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$32"
$wb.Names.Item("web").RefersTo = "='#system'!`$AA`$2:`$AA`$144"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$AB`$2:`$AB`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AC`$2:`$AC`$10"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AD`$2:`$AD`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AE`$2:`$AE`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AF`$2:`$AF`$27"
$wb.Names.Add("tn.5250", "='#system'!`$Z`$2:`$Z`$6")
